$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("General")
$ws2 = $wb.Worksheets.Item("Gana")

# --- "Gana" sheet: add new row 14 (BDD-ified "upload resume" step) ---
$ws2.Range("A14").Value = " Fill out the General Information Form"
$ws2.Range("A14").WrapText = $true
$ws2.Range("A14").VerticalAlignment = -4108

$ws2.Range("B14").Value = 13
$ws2.Range("B14").HorizontalAlignment = 1
$ws2.Range("B14").WrapText = $true
$ws2.Range("B14").VerticalAlignment = -4108

$ws2.Range("C14").Value = ' Upload "Resume" in the specified format.'

$ws2.Range("D14").Value = "//input[@type='file']"
$ws2.Range("D14").WrapText = $true
$ws2.Range("D14").VerticalAlignment = -4108
$ws2.Range("D14").Font.Name = "Arial Unicode MS"
$ws2.Range("D14").Font.Size = 10

$ws2.Range("E14").WrapText = $true
$ws2.Range("E14").VerticalAlignment = -4108
$ws2.Range("E14").Font.Name = "Arial Unicode MS"
$ws2.Range("E14").Font.Size = 10

$ws2.Range("F14").HorizontalAlignment = -4131
$ws2.Range("F14").VerticalAlignment = -4160
$ws2.Range("F14").WrapText = $true

# --- view state: drop the scrolled topLeftCell and move the selection ---
$ws1.Activate()
$ws1.Range("A7:XFD7").Select()

$ws2.Activate()
$ws2.Range("C14").Select()
